$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available ."
$ws.Range("C2").Value = "'When exceeding nnn, Obstacle Avoidance is not available"
$ws.Range("D2").Value = "'7-14"
$ws.Range("E2").Value = "'Missing"

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available ."
$ws.Range("C3").Value = "'When exceeding nnn,"
$ws.Range("D3").Value = "'7-9"
$ws.Range("E3").Value = "'False"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available ."
$ws.Range("C4").Value = "'Obstacle Avoidance is not available"
$ws.Range("D4").Value = "'10-14"
$ws.Range("E4").Value = "'False"

$ws.Range("A5").Value = 14
$ws.Range("B5").Value = "'Aircraft in Restricted Zones . Check map to find Recommended Zones ."
$ws.Range("C5").Value = "'Aircraft in Restricted Zones"
$ws.Range("D5").Value = "'0-3"
$ws.Range("E5").Value = "'Missing"

$ws.Range("A6").Value = 15
$ws.Range("B6").Value = "'Aircraft is close to the Home Point . Initiating Return to Home will now trigger Auto Landing ."
$ws.Range("C6").Value = "'Aircraft is close to the Home Point"
$ws.Range("D6").Value = "'0-6"
$ws.Range("E6").Value = "'Missing"

$ws.Range("A7").Value = 19
$ws.Range("B7").Value = "'Aircraft is tilted , please keep the aircraft stationary and level before flight ."
$ws.Range("C7").Value = "'Aircraft is tilted"
$ws.Range("D7").Value = "'0-2"
$ws.Range("E7").Value = "'Missing"

$ws.Range("A8").Value = 19
$ws.Range("B8").Value = "'Aircraft is tilted , please keep the aircraft stationary and level before flight ."
$ws.Range("C8").Value = "'please keep the aircraft stationary and level before flight"
$ws.Range("D8").Value = "'4-12"
$ws.Range("E8").Value = "'Missing"

$ws.Range("A9").Value = 19
$ws.Range("B9").Value = "'Aircraft is tilted , please keep the aircraft stationary and level before flight ."
$ws.Range("C9").Value = "'Aircraft is tilted , please keep the aircraft stationary and level before flight"
$ws.Range("D9").Value = "'0-12"
$ws.Range("E9").Value = "'False"

$ws.Range("A10").Value = 29
$ws.Range("B10").Value = "'Another aircraft is dangerously close , please descend to a safer altitude ."
$ws.Range("C10").Value = "'Another aircraft is dangerously close"
$ws.Range("D10").Value = "'0-4"
$ws.Range("E10").Value = "'Missing"

$ws.Range("A11").Value = 42
$ws.Range("B11").Value = "'Camera sensor error . Hardware malfunction : Contact DJI Support to arrange for repairs ."
$ws.Range("C11").Value = "'Hardware malfunction"
$ws.Range("D11").Value = "'4-5"
$ws.Range("E11").Value = "'Missing"

$ws.Range("A12").Value = 50
$ws.Range("B12").Value = "'Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C12").Value = "'2. Calibrate Compass Before Takeoff"
$ws.Range("D12").Value = "'35-39"
$ws.Range("E12").Value = "'Missing"

$ws.Range("A13").Value = 50
$ws.Range("B13").Value = "'Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C13").Value = "'2."
$ws.Range("D13").Value = "'35-35"
$ws.Range("E13").Value = "'False"

$ws.Range("A14").Value = 50
$ws.Range("B14").Value = "'Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C14").Value = "'Calibrate Compass Before Takeoff"
$ws.Range("D14").Value = "'36-39"
$ws.Range("E14").Value = "'False"

$ws.Range("A15").Value = 66
$ws.Range("B15").Value = "'Downlink data connection lost for nnn seconds ."
$ws.Range("C15").Value = "'Downlink data connection lost for nnn seconds"
$ws.Range("D15").Value = "'0-6"
$ws.Range("E15").Value = "'Missing"

$ws.Range("A16").Value = 66
$ws.Range("B16").Value = "'Downlink data connection lost for nnn seconds ."
$ws.Range("C16").Value = "'Downlink data connection lost for nnn"
$ws.Range("D16").Value = "'0-5"
$ws.Range("E16").Value = "'False"

$ws.Range("A17").Value = 81
$ws.Range("B17").Value = "'Extra payload detected . Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety ."
$ws.Range("C17").Value = "'Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety"
$ws.Range("D17").Value = "'4-23"
$ws.Range("E17").Value = "'Missing"

$ws.Range("A18").Value = 91
$ws.Range("B18").Value = "'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C18").Value = "'GEO Zone Info: The target area is in an Altitude Zone"
$ws.Range("D18").Value = "'0-10"
$ws.Range("E18").Value = "'Missing"

$ws.Range("A19").Value = 91
$ws.Range("B19").Value = "'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C19").Value = "'GEO Zone Info:"
$ws.Range("D19").Value = "'0-2"
$ws.Range("E19").Value = "'False"

$ws.Range("A20").Value = 91
$ws.Range("B20").Value = "'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C20").Value = "'The target area is in an Altitude Zone"
$ws.Range("D20").Value = "'3-10"
$ws.Range("E20").Value = "'False"

$ws.Range("A21").Value = 99
$ws.Range("B21").Value = "'GPS signal low . Aircraft unable to auto hover and takeoff restricted . Move to environment with adequate light . Unlocking takeoff restrictions not recommended ."
$ws.Range("C21").Value = "'Unlocking takeoff restrictions not recommended"
$ws.Range("D21").Value = "'20-24"
$ws.Range("E21").Value = "'Missing"

$ws.Range("A22").Value = 104
$ws.Range("B22").Value = "'Home Point Recorded , Return-to-Home Altitude : 98FT ."
$ws.Range("C22").Value = "'Home Point Recorded"
$ws.Range("D22").Value = "'0-2"
$ws.Range("E22").Value = "'Missing"

$ws.Range("A23").Value = 113
$ws.Range("B23").Value = "'Insufficient SD card space . Change card or delete files ."
$ws.Range("C23").Value = "'Insufficient SD card space"
$ws.Range("D23").Value = "'0-3"
$ws.Range("E23").Value = "'Missing"

$ws.Range("A24").Value = 123
$ws.Range("B24").Value = "'Motor Obstructed . Propulsion output is limited to ensure the health of the battery ."
$ws.Range("C24").Value = "'Propulsion output is limited to ensure the health of the battery"
$ws.Range("D24").Value = "'3-13"
$ws.Range("E24").Value = "'Missing"

$ws.Range("A25").Value = 125
$ws.Range("B25").Value = "'No GPS signal . Unable to hover . Fly with caution ."
$ws.Range("C25").Value = "'Unable to hover"
$ws.Range("D25").Value = "'4-6"
$ws.Range("E25").Value = "'Missing"

$ws.Range("A26:E32").EntireRow.Delete()
